$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data > Sort: the whole table (header row 1, data A2:G65) sorted
#     ascending by "Jumlah Kombinasi" (col A) then "Kombinasi" (col B) ---
$fullRange = $ws.Range("A1:G65")
$key1 = $ws.Range("A1")
$key2 = $ws.Range("B1")
$fullRange.Sort($key1, 1, $key2, $null, 1, $null, 1, 1, $false, 1, 1, 1) | Out-Null

# --- A second, narrower sort was also performed on the current selection
#     (rows 2-7, already in order, so it leaves the data unchanged but
#     registers its own sort state) ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$subKey1 = $ws.Range("A2:A7")
$subKey2 = $ws.Range("B2:B7")
$sortObj.SortFields.Add($subKey1) | Out-Null
$sortObj.SortFields.Add($subKey2) | Out-Null
$sortObj.SetRange($ws.Range("A2:G7")) | Out-Null
$sortObj.Header = 2
$sortObj.Apply() | Out-Null

# --- Final selection left on the sheet ---
$ws.Range("B2:F8").Select() | Out-Null
